$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need an explicit Text
# number format first, otherwise the COM Value setter auto-coerces the
# string into a numeric value (losing the original text formatting).
$ws.Range("D2").Value = "66.146.48"
$ws.Range("E2").Value = "  +5.97%  "
$ws.Range("D3").Value = "2.998.48"
$ws.Range("E3").Value = "  +3.57%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.27"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.86"
$ws.Range("E6").Value = "  +13.47%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +3.59%  "
$ws.Range("D9").Value = "2.995.57"
$ws.Range("E9").Value = "  +3.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.61"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("E11").Value = "  +3.62%  "
$ws.Range("E12").Value = "  +5.68%  "
$ws.Range("E13").Value = "  +5.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.53"
$ws.Range("E14").Value = "  +5.65%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "66.135.00"
$ws.Range("E16").Value = "  +6.11%  "
$ws.Range("D17").Value = "3.494.53"
$ws.Range("E17").Value = "  +3.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.90"
$ws.Range("E18").Value = "  +5.89%  "
$ws.Range("D19").Value = "2.998.90"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "452.39"
$ws.Range("E20").Value = "  +6.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.84"
$ws.Range("E21").Value = "  +6.43%  "
$ws.Range("E22").Value = "  +4.45%  "
$ws.Range("E23").Value = "  +7.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.32"
$ws.Range("E24").Value = "  +4.84%  "
$ws.Range("E25").Value = "  +14.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.26"
$ws.Range("E26").Value = "  +3.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.30"
$ws.Range("E27").Value = "  +2.83%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +14.82%  "
$ws.Range("E30").Value = "  +20.17%  "
$ws.Range("E31").Value = "  +5.97%  "
$ws.Range("E32").Value = "  -5.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.18"
$ws.Range("E33").Value = "  +5.88%  "
$ws.Range("E34").Value = "  +5.27%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.992"
$ws.Range("E36").Value = "  +4.62%  "
$ws.Range("E37").Value = "  +8.11%  "
$ws.Range("E38").Value = "  +8.99%  "
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("E41").Value = "  +16.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "44.11"
$ws.Range("E42").Value = "  +7.61%  "
$ws.Range("E43").Value = "  +7.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.44"
$ws.Range("E44").Value = "  +5.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "400.72"
$ws.Range("E45").Value = "  +13.31%  "
$ws.Range("E46").Value = "  +6.48%  "
$ws.Range("D47").Value = "2.767.08"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.41"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.79"
$ws.Range("E50").Value = "  +12.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.106"
$ws.Range("E51").Value = "  +4.02%  "
